$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: every cell below is stored as plain text in the source data (prices
# like "1.00" / "0.581" and percent strings like "  +3.66%  "). Forcing the
# cell to Text number-format before assigning keeps Excel from silently
# coercing these into numbers (which would drop trailing zeros / exact
# decimal text) - matches the original file's text-cell semantics.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.661.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.496.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.12'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.63'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.80%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.495.44'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.44%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.70%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.89%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.103.88'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.08%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.28'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.56%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.673.75'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.58%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.487.93'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.08'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.73'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.62%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.17%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.99%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.39'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.182'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.978'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.35'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.08%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.71%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.56'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.92%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.61'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.20%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.50%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.35%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.70%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.81'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.83%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0747'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.46'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.822.93'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.67'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.01'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '354.89'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.42%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.76'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +11.14%  '
